$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (trailing / leading spaces are intentional)
$ws.Range("A1").Value = "Description "
$ws.Range("B1").Value = " Quantité"
$ws.Range("C1").Value = "Prix unitaire (€) "
$ws.Range("D1").Value = "TVA (%) "

# Widen column C slightly to fit the new header text
$ws.Columns.Item(3).ColumnWidth = 15.5

# New trailing row holding a single blank/space value in column A
$ws.Range("A9").Value = " "

# Turn the data range into a real table (ListObject) so imports can
# gracefully skip/ignore missing information
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:D9"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium24"
